$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices + 1h volume %) -- GitHub Actions refresh
# Force the Price column (D) to stay text so values like trailing zeros survive,
# matching the original inlineStr cell type for every row touched below.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.565.74'
$ws.Range("E2").Value = '  -2.53%  '
$ws.Range("D3").Value = '1.656.24'
$ws.Range("E3").Value = '  -4.19%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '213.94'
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '23.96'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("D10").Value = '0.0619'
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D12").Value = '1.888.74'
$ws.Range("E12").Value = '  -4.46%  '
$ws.Range("D13").Value = '1.664.22'
$ws.Range("E13").Value = '  -3.89%  '
$ws.Range("D14").Value = '4.14'
$ws.Range("E14").Value = '  -2.45%  '
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '65.80'
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").Value = '27.533.04'
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("D18").Value = '240.52'
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("E19").Value = '  -3.14%  '
$ws.Range("D20").Value = '7.55'
$ws.Range("E20").Value = '  -4.39%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '4.47'
$ws.Range("E22").Value = '  -4.01%  '
$ws.Range("D23").Value = '9.32'
$ws.Range("E23").Value = '  -3.73%  '
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").Value = '145.71'
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("E26").Value = '  -3.95%  '
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.0501'
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.21'
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = '3.32'
$ws.Range("E32").Value = '  -2.91%  '
$ws.Range("D33").Value = '1.448.75'
$ws.Range("E33").Value = '  -2.74%  '
$ws.Range("E34").Value = '  -5.24%  '
$ws.Range("E35").Value = '  -5.44%  '
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("D37").Value = '0.921'
$ws.Range("E37").Value = '  -5.96%  '
$ws.Range("E38").Value = '  -3.05%  '
$ws.Range("E39").Value = '  -5.43%  '
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  -2.78%  '
$ws.Range("D41").Value = '69.06'
$ws.Range("E41").Value = '  -1.33%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '5.40'
$ws.Range("E43").Value = '  -4.61%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D46").Value = '1.796.88'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = '88.42'
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("E49").Value = '  -6.36%  '
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").Value = '7.80'
$ws.Range("E51").Value = '  -4.80%  '
